$wb = $excel.ActiveWorkbook

# --- Sheet: Folder Inventory ---
$ws = $wb.Worksheets.Item("Folder Inventory")

# Insert a new row at the top of the data (row 2), pushing everything else down.
$ws.Rows.Item(2).Insert()
# Excel's insert picks up formatting from the row above (the bold header row);
# clear that back out so the new row matches the plain data rows.
$ws.Rows.Item(2).ClearFormats()

$newTitle = "Work with Data Lake and Data Factory Pipelines in Microsoft Fabric" + [char]0x200B
$ws.Cells.Item(2, 1).Value = $newTitle
$ws.Cells.Item(2, 2).Value = $newTitle
$ws.Cells.Item(2, 3).Value = "2025-06-12 17:26:19 +0530"
$ws.Cells.Item(2, 4).Value = 1
$ws.Cells.Item(2, 5).Value = "Root"

# --- Sheet: Metadata ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Cells.Item(3, 2).Value = "2025-06-12 11:56:38 UTC"
$wsMeta.Cells.Item(4, 2).Value = 73
# "Workflow Run" is stored as text ("12"), not a number; force text storage
# then drop the quote-prefix style Excel adds so no stray style index sticks.
$wsMeta.Cells.Item(5, 2).NumberFormat = "@"
$wsMeta.Cells.Item(5, 2).Value = "12"
$wsMeta.Cells.Item(5, 2).ClearFormats()

# --- Sheet: Summary ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Cells.Item(2, 2).Value = 73
$wsSummary.Cells.Item(3, 2).Value = 73
$wsSummary.Cells.Item(5, 2).Value = "2025-06-12 17:26:19 +0530"
